$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "68.153.94"
Set-TextValue "E2" "  -0.51%  "
Set-TextValue "D3" "3.271.78"
Set-TextValue "E3" "  +0.58%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "584.08"
Set-TextValue "E5" "  +0.00%  "
Set-TextValue "D6" "184.78"
Set-TextValue "E6" "  +1.98%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "0.602"
Set-TextValue "E8" "  +1.23%  "
Set-TextValue "E9" "  -2.86%  "
Set-TextValue "E10" "  -0.75%  "
Set-TextValue "E11" "  -2.39%  "
Set-TextValue "D12" "3.843.35"
Set-TextValue "E12" "  +0.71%  "
Set-TextValue "E13" "  +0.92%  "
Set-TextValue "D14" "27.58"
Set-TextValue "E14" "  -2.26%  "
Set-TextValue "D15" "68.182.05"
Set-TextValue "E15" "  -0.47%  "
Set-TextValue "D16" "0.0000168"
Set-TextValue "E16" "  -1.52%  "
Set-TextValue "D17" "3.273.41"
Set-TextValue "E17" "  +1.03%  "
Set-TextValue "D18" "5.75"
Set-TextValue "E18" "  -1.08%  "
Set-TextValue "D19" "13.36"
Set-TextValue "E19" "  -0.86%  "
Set-TextValue "D20" "417.90"
Set-TextValue "E20" "  +6.04%  "
Set-TextValue "D21" "7.57"
Set-TextValue "E21" "  -1.15%  "
Set-TextValue "B22" "Dai"
Set-TextValue "C22" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D22" "1.00"
Set-TextValue "E22" "  +0.40%  "
Set-TextValue "B23" "Litecoin"
Set-TextValue "C23" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D23" "71.45"
Set-TextValue "E23" "  -0.26%  "
Set-TextValue "E24" "  -1.33%  "
Set-TextValue "D25" "0.0000117"
Set-TextValue "E25" "  -0.88%  "
Set-TextValue "E26" "  -1.26%  "
Set-TextValue "D27" "9.45"
Set-TextValue "E27" "  -1.54%  "
Set-TextValue "D28" "0.998"
Set-TextValue "E28" "  -0.19%  "
Set-TextValue "D29" "1.95"
Set-TextValue "E29" "  -1.72%  "
Set-TextValue "E30" "  -0.94%  "
Set-TextValue "D31" "5.48"
Set-TextValue "E31" "  -3.39%  "
Set-TextValue "D32" "6.89"
Set-TextValue "E32" "  -2.99%  "
Set-TextValue "E33" "  +0.02%  "
Set-TextValue "E34" "  -1.86%  "
Set-TextValue "D35" "163.92"
Set-TextValue "E35" "  -0.28%  "
Set-TextValue "E36" "  -2.49%  "
Set-TextValue "D37" "1.89"
Set-TextValue "E37" "  -1.95%  "
Set-TextValue "D38" "27.15"
Set-TextValue "E38" "  +3.42%  "
Set-TextValue "E39" "  -2.75%  "
Set-TextValue "D40" "4.47"
Set-TextValue "E40" "  -2.57%  "
Set-TextValue "D41" "6.34"
Set-TextValue "E41" "  -3.66%  "
Set-TextValue "D42" "2.665.85"
Set-TextValue "E42" "  +3.07%  "
Set-TextValue "D43" "40.89"
Set-TextValue "E43" "  -1.35%  "
Set-TextValue "E44" "  -1.05%  "
Set-TextValue "E45" "  -1.55%  "
Set-TextValue "D46" "337.73"
Set-TextValue "E46" "  -1.37%  "
Set-TextValue "D47" "24.51"
Set-TextValue "E47" "  -0.53%  "
Set-TextValue "E48" "  -2.41%  "
Set-TextValue "E49" "  -0.13%  "
Set-TextValue "D50" "0.980"
Set-TextValue "E50" "  +0.20%  "
Set-TextValue "E51" "  -1.19%  "
